$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 12262
$ws.Range("D2").Value = 151
$ws.Range("E2").Value = 0.593
$ws.Range("F2").Value = 29.7

# Row 3
$ws.Range("C3").Value = 8979
$ws.Range("D3").Value = 188
$ws.Range("F3").Value = 31.2

# Row 4
$ws.Range("C4").Value = 29466
$ws.Range("E4").Value = 0.572
$ws.Range("F4").Value = 18.5

# Row 5
$ws.Range("C5").Value = 42975
$ws.Range("D5").Value = 444
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 42.9
